$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("n6", "n6_e280_4321_2_7_UnitedStoneInternational_Cleveland_OH_.jpeg", "True", "no_meltpatch"),
    @("n7", "n7_e274_4321_1_2_UnitedStoneInternational_Cleveland_OH_.jpeg", "True", "no_meltpatch"),
    @("n8", "n8_e276_4321_1_4_UnitedStoneInternational_Cleveland_OH_.jpeg", "True", "no_meltpatch"),
    @("n9", "n9_e277_4321_1_5_UnitedStoneInternational_Cleveland_OH_.jpeg", "True", "no_meltpatch"),
    @("n10", "n10_e279_4321_1_7_UnitedStoneInternational_Cleveland_OH_.jpeg", "True", "no_meltpatch")
)

$startRow = 7
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = "'" + $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
